# Add a new day (column H) to the time-tracking sheet, mirroring the
# data already recorded for the previous day in column G:
#   - Row 1: header date formula (=previous day + 1)
#   - Row 3: the recurring "component" formula (=$A32)
#   - Rows 4-10: clock-in/out / break time-of-day values
#   - Row 21: "Finished at" time-of-day value
# Formatting for each new cell is copied from the corresponding G-column
# cell (via PasteSpecial formats) so the same shared cell style is reused
# rather than creating new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header date - one day after G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Formula = "=G1+1"

# Row 3: recurring "component" formula, same as the rest of the row
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Formula = "=`$A32"

# Rows 4-10: time-of-day values for the new day
$ws.Range("G4").Copy()
$ws.Range("H4:H10").PasteSpecial(-4122)
$ws.Range("H4").Value = 0.45833333333333331
$ws.Range("H5").Value = 0.47916666666666669
$ws.Range("H6").Value = 0.5
$ws.Range("H7").Value = 0.52083333333333337
$ws.Range("H8").Value = 0.58333333333333337
$ws.Range("H9").Value = 0.64930555555555558
$ws.Range("H10").Value = 0.79999999999999993

# Row 21: "Finished at" time for the new day
$ws.Range("G21").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H21").Value = 0.83124999999999993

$ws.Application.CutCopyMode = $false

# Move the selection to H21, matching the saved view state
$ws.Range("H21").Select()
